# Auto-generated Excel COM-interop script applying the Tiamat_Profits.xlsx diff.
# Updates cached market-price / profit figures across the ALC, ARM, BSM, CRP,
# GSM, LTW and WVR sheets (columns H-N) to refreshed values from the latest
# scheduled market-data run. A couple of rows also gain/lose an 'M' (LeveProfitNQ)
# cell where the NQ price feed newly has/lost a data point.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 53003
$ws.Range("J13").Value = 49603.6
$ws.Range("L13").Value = 49603.6
$ws.Range("N13").Value = -49941.6

$ws.Range("H38").Value = 468.09525
$ws.Range("I38").Value = 88.916664
$ws.Range("J38").Value = 973.6667
$ws.Range("K38").Value = 266.749992
$ws.Range("L38").Value = 2921.0001
$ws.Range("M38").Value = 105.250008
$ws.Range("N38").Value = -3665.0001

$ws.Range("H39").Value = 1330.0476
$ws.Range("I39").Value = 119
$ws.Range("J39").Value = 1531.8889
$ws.Range("K39").Value = 357
$ws.Range("L39").Value = 4595.6667
$ws.Range("M39").Value = -61
$ws.Range("N39").Value = -5187.6667

$ws.Range("H69").Value = 4904094
$ws.Range("J69").Value = 3200
$ws.Range("L69").Value = 9600
$ws.Range("N69").Value = -11348

$ws.Range("H72").Value = 4904094
$ws.Range("J72").Value = 3200
$ws.Range("L72").Value = 28800
$ws.Range("N72").Value = -37536

$ws.Range("H132").Value = 191573.3
$ws.Range("I132").Value = 3473.1714
$ws.Range("J132").Value = 557323.5600000001
$ws.Range("K132").Value = 10419.5142
$ws.Range("L132").Value = 1671970.68
$ws.Range("M132").Value = -7889.514200000001
$ws.Range("N132").Value = -1677030.68

$ws.Range("H138").Value = 1237.0133
$ws.Range("I138").Value = 605.6889
$ws.Range("J138").Value = 2184
$ws.Range("K138").Value = 1817.0667
$ws.Range("L138").Value = 6552
$ws.Range("M138").Value = 3322.9333
$ws.Range("N138").Value = -16832


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 71256.75
$ws.Range("J19").Value = 71256.75
$ws.Range("L19").Value = 71256.75
$ws.Range("N19").Value = -71714.75

$ws.Range("H33").Value = 58021.75
$ws.Range("I33").Value = 22000
$ws.Range("K33").Value = 22000
$ws.Range("M33").Value = -21671

$ws.Range("H110").Value = 728.25
$ws.Range("I110").Value = 666.5789
$ws.Range("J110").Value = 962.6
$ws.Range("K110").Value = 666.5789
$ws.Range("L110").Value = 962.6
$ws.Range("M110").Value = 1378.4211
$ws.Range("N110").Value = -5052.6


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 55004.5
$ws.Range("J14").Value = 55004.5
$ws.Range("L14").Value = 55004.5
$ws.Range("N14").Value = -55348.5

$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

$ws.Range("H105").Value = 996262.6
$ws.Range("I105").Value = 1991250.2
$ws.Range("J105").Value = 1275
$ws.Range("K105").Value = 1991250.2
$ws.Range("L105").Value = 1275
$ws.Range("M105").Value = -1989503.2
$ws.Range("N105").Value = -4769


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

$ws.Range("H58").Value = 2180.05
$ws.Range("I58").Value = 717.4783
$ws.Range("J58").Value = 4158.8237
$ws.Range("K58").Value = 717.4783
$ws.Range("L58").Value = 4158.8237
$ws.Range("M58").Value = -514.4783
$ws.Range("N58").Value = -4564.8237

$ws.Range("H98").Value = 54890
$ws.Range("J98").Value = 54890
$ws.Range("L98").Value = 54890
$ws.Range("N98").Value = -59382

$ws.Range("H99").Value = 4883.3335
$ws.Range("I99").Value = 4425
$ws.Range("J99").Value = 5800
$ws.Range("K99").Value = 4425
$ws.Range("L99").Value = 5800
$ws.Range("M99").Value = -2927
$ws.Range("N99").Value = -8796

$ws.Range("H106").Value = 59993.332
$ws.Range("J106").Value = 59993.332
$ws.Range("L106").Value = 59993.332
$ws.Range("N106").Value = -62517.332

$ws.Range("H126").Value = 4883.3335
$ws.Range("I126").Value = 4425
$ws.Range("J126").Value = 5800
$ws.Range("K126").Value = 13275
$ws.Range("L126").Value = 17400
$ws.Range("M126").Value = -10805
$ws.Range("N126").Value = -22340

$ws.Range("H136").Value = 2180.05
$ws.Range("I136").Value = 717.4783
$ws.Range("J136").Value = 4158.8237
$ws.Range("K136").Value = 2152.4349
$ws.Range("L136").Value = 12476.4711
$ws.Range("M136").Value = 397.5650999999998
$ws.Range("N136").Value = -17576.4711


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 43754.5
$ws.Range("J25").Value = 55006
$ws.Range("L25").Value = 55006
$ws.Range("N25").Value = -56064

$ws.Range("H113").Value = 665
$ws.Range("I113").Value = 666.6667
$ws.Range("J113").Value = 664.2857
$ws.Range("K113").Value = 666.6667
$ws.Range("L113").Value = 664.2857
$ws.Range("M113").Value = 1503.3333
$ws.Range("N113").Value = -5004.2857

$ws.Range("H126").Value = 2584
$ws.Range("J126").Value = 2480
$ws.Range("L126").Value = 7440
$ws.Range("N126").Value = -12380

$ws.Range("H129").Value = 46955.4
$ws.Range("J129").Value = 46955.4
$ws.Range("L129").Value = 46955.4
$ws.Range("N129").Value = -56955.4


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 58340
$ws.Range("J5").Value = 58340
$ws.Range("L5").Value = 58340
$ws.Range("N5").Value = -58566

$ws.Range("H7").Value = 6107.273
$ws.Range("I7").Value = 8600
$ws.Range("J7").Value = 4030
$ws.Range("K7").Value = 8600
$ws.Range("L7").Value = 4030
$ws.Range("M7").Value = -8488
$ws.Range("N7").Value = -4254

$ws.Range("H22").Value = 1989
$ws.Range("I22").Value = 900
$ws.Range("J22").Value = 2805.75
$ws.Range("K22").Value = 900
$ws.Range("L22").Value = 2805.75
$ws.Range("M22").Value = -605
$ws.Range("N22").Value = -3395.75

$ws.Range("H27").Value = 1989
$ws.Range("I27").Value = 900
$ws.Range("J27").Value = 2805.75
$ws.Range("K27").Value = 900
$ws.Range("L27").Value = 2805.75
$ws.Range("M27").Value = -793
$ws.Range("N27").Value = -3019.75

$ws.Range("H68").Value = 2334.2
$ws.Range("I68").Value = 1799.1111
$ws.Range("J68").Value = 3136.8333
$ws.Range("K68").Value = 1799.1111
$ws.Range("L68").Value = 3136.8333
$ws.Range("M68").Value = -1050.1111
$ws.Range("N68").Value = -4634.8333

$ws.Range("H71").Value = 2334.2
$ws.Range("I71").Value = 1799.1111
$ws.Range("J71").Value = 3136.8333
$ws.Range("K71").Value = 8995.5555
$ws.Range("L71").Value = 15684.1665
$ws.Range("M71").Value = -5251.5555
$ws.Range("N71").Value = -23172.1665

$ws.Range("H94").Value = 30765.8
$ws.Range("J94").Value = 30765.8
$ws.Range("L94").Value = 30765.8
$ws.Range("N94").Value = -32117.8

$ws.Range("H126").Value = 6107.273
$ws.Range("I126").Value = 8600
$ws.Range("J126").Value = 4030
$ws.Range("K126").Value = 25800
$ws.Range("L126").Value = 12090
$ws.Range("M126").Value = -23330
$ws.Range("N126").Value = -17030


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 65012.75
$ws.Range("I21").Value = 50000
$ws.Range("K21").Value = 50000
$ws.Range("M21").Value = -49765

$ws.Range("H35").Value = 65012.75
$ws.Range("I35").Value = 50000
$ws.Range("K35").Value = 50000
$ws.Range("M35").Value = -49710

$ws.Range("H126").Value = 1068
$ws.Range("I126").Value = 960
$ws.Range("K126").Value = 2880
$ws.Range("M126").Value = -410
